$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing rows that are no longer part of the n=45 dataset
$ws.Rows("47:51").Delete()

# Updated Composite_Reward / Composite_Reward_Squared values for n = 45 subjects
$data = @(
    @(2, 1001, -0.62222222222222179, 0.38716049382715995),
    @(3, 1006, -0.62222222222222179, 0.38716049382715995),
    @(4, 1009, 0.37777777777777821, 0.14271604938271637),
    @(5, 1010, -1.6222222222222218, 2.6316049382716034),
    @(6, 1011, 2.3777777777777782, 5.6538271604938295),
    @(7, 1012, -4.6222222222222218, 21.364938271604935),
    @(8, 1013, 4.3777777777777782, 19.164938271604942),
    @(9, 1015, -3.6222222222222218, 13.120493827160491),
    @(10, 1016, 2.3777777777777782, 5.6538271604938295),
    @(11, 1019, 2.3777777777777782, 5.6538271604938295),
    @(12, 1021, -1.6222222222222218, 2.6316049382716034),
    @(13, 1242, 2.3777777777777782, 5.6538271604938295),
    @(14, 1243, -2.6222222222222218, 6.8760493827160474),
    @(15, 1244, 4.3777777777777782, 19.164938271604942),
    @(16, 1248, 2.3777777777777782, 5.6538271604938295),
    @(17, 1249, 1.3777777777777782, 1.8982716049382728),
    @(18, 1251, -4.6222222222222218, 21.364938271604935),
    @(19, 1255, -1.6222222222222218, 2.6316049382716034),
    @(20, 1276, -2.6222222222222218, 6.8760493827160474),
    @(21, 1282, 4.3777777777777782, 19.164938271604942),
    @(22, 1286, -0.62222222222222179, 0.38716049382715995),
    @(23, 1294, -2.6222222222222218, 6.8760493827160474),
    @(24, 1301, -3.6222222222222218, 13.120493827160491),
    @(25, 1302, 4.3777777777777782, 19.164938271604942),
    @(26, 1303, 0.37777777777777821, 0.14271604938271637),
    @(27, 3116, 1.3777777777777782, 1.8982716049382728),
    @(28, 3122, 3.3777777777777782, 11.409382716049386),
    @(29, 3125, 4.3777777777777782, 19.164938271604942),
    @(30, 3140, -1.6222222222222218, 2.6316049382716034),
    @(31, 3143, 0.37777777777777821, 0.14271604938271637),
    @(32, 3152, 4.3777777777777782, 19.164938271604942),
    @(33, 3166, -3.6222222222222218, 13.120493827160491),
    @(34, 3167, -2.6222222222222218, 6.8760493827160474),
    @(35, 3170, -4.6222222222222218, 21.364938271604935),
    @(36, 3173, 3.3777777777777782, 11.409382716049386),
    @(37, 3176, -3.6222222222222218, 13.120493827160491),
    @(38, 3189, -0.62222222222222179, 0.38716049382715995),
    @(39, 3190, 3.3777777777777782, 11.409382716049386),
    @(40, 3199, 0.37777777777777821, 0.14271604938271637),
    @(41, 3200, -3.6222222222222218, 13.120493827160491),
    @(42, 3206, -1.6222222222222218, 2.6316049382716034),
    @(43, 3210, 2.3777777777777782, 5.6538271604938295),
    @(44, 3212, -3.6222222222222218, 13.120493827160491),
    @(45, 3218, 0.37777777777777821, 0.14271604938271637),
    @(46, 3220, 1.3777777777777782, 1.8982716049382728)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
